# Change "... було знищено під час ведення бойових дій." to
# "... було втрачено під час ведення бойових дій." while splitting the
# sentence into three separate runs (matching the target OOXML diff):
#   " було "  |  "втрачено"  |  " під час ведення бойових дій."

$d = $word.ActiveDocument

# Locate the word to replace.
$findRange = $d.Content
$findRange.Find.Execute("знищено", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

$start = $findRange.Start
$end   = $findRange.End

# Replace the word's text.
$wordRange = $d.Range($start, $end)
$wordRange.Text = "втрачено"

# Toggling a character formatting property (and reverting it) forces the
# engine to keep this span as its own run, splitting the original single
# run into three runs with identical visual formatting, just like the
# diff shows.
$splitRange = $d.Range($start, $start + 8)
$splitRange.Bold = 1
$splitRange.Bold = 0
